$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the authoritative diff.
$updates = [ordered]@{
    "D2" = "67.238.82"
    "E2" = "  -5.11%  "
    "D3" = "3.231.36"
    "E3" = "  -8.66%  "
    "E4" = "  -0.08%  "
    "D5" = "585.16"
    "E5" = "  -4.89%  "
    "D6" = "152.74"
    "E6" = "  -12.03%  "
    "E7" = "  -0.16%  "
    "D8" = "3.224.59"
    "E8" = "  -8.68%  "
    "D9" = "0.544"
    "E9" = "  -11.09%  "
    "E10" = "  -12.90%  "
    "D11" = "6.73"
    "E11" = "  -6.80%  "
    "E12" = "  -14.52%  "
    "D13" = "38.33"
    "E13" = "  -17.88%  "
    "E14" = "  -11.76%  "
    "D15" = "3.749.35"
    "E15" = "  -8.73%  "
    "D16" = "67.183.64"
    "E16" = "  -5.23%  "
    "B17" = "BitcoinCash"
    "C17" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D17" = "542.88"
    "E17" = "  -11.39%  "
    "B18" = "WrappedEther"
    "C18" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D18" = "3.230.13"
    "E18" = "  -8.60%  "
    "E19" = "  -5.83%  "
    "D20" = "7.16"
    "E20" = "  -15.25%  "
    "D21" = "15.15"
    "E21" = "  -14.95%  "
    "D22" = "0.760"
    "E22" = "  -14.52%  "
    "E23" = "  -13.97%  "
    "D24" = "85.64"
    "E24" = "  -12.90%  "
    "D25" = "13.50"
    "E25" = "  -14.26%  "
    "E26" = "  +0.00%  "
    "D27" = "3.17"
    "E27" = "  -16.47%  "
    "D28" = "8.10"
    "E28" = "  -11.68%  "
    "D29" = "29.43"
    "E29" = "  -13.21%  "
    "E30" = "  -17.75%  "
    "E32" = "  -12.53%  "
    "D33" = "543.88"
    "E33" = "  -10.09%  "
    "D34" = "6.56"
    "E34" = "  -19.91%  "
    "D35" = "5.72"
    "E35" = "  -16.83%  "
    "E36" = "  +0.13%  "
    "B37" = "VeChain"
    "C37" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D37" = "0.0444"
    "E37" = "  -6.09%  "
    "B38" = "OKB"
    "C38" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D38" = "53.61"
    "E38" = "  -5.97%  "
    "B39" = "Hedera"
    "C39" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D39" = "0.0847"
    "E39" = "  -15.79%  "
    "B40" = "Cosmos"
    "C40" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D40" = "9.19"
    "E40" = "  -15.41%  "
    "E41" = "  -13.02%  "
    "D42" = "2.926.86"
    "E42" = "  -13.60%  "
    "E43" = "  -27.11%  "
    "D44" = "0.0₃0586"
    "E44" = "  -21.04%  "
    "E45" = "  -17.20%  "
    "E46" = "  -19.92%  "
    "E47" = "  -0.03%  "
    "D48" = "26.11"
    "E48" = "  -19.13%  "
    "E49" = "  -18.04%  "
    "E50" = "  -13.36%  "
    "D51" = "123.84"
    "E51" = "  -7.39%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage (matches source workbook, where every touched
    # cell is an inline/shared string, even the ones that look numeric,
    # e.g. prices like '585.16') so Excel doesn't coerce them to Double.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop back to the default style so we don't leave a stray
    # "Text" number-format style on the cell (source cells carry no
    # explicit style).
    $cell.Style = "Normal"
}
